$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.05980000000000001
$ws.Range("E2").Value = -0.0365
$ws.Range("F2").Value = 0.15
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 81.9
$ws.Range("L2").Value = 0.338010730499381
$ws.Range("M2").Value = 51.461
$ws.Range("N2").Value = 0.04234079315451703
$ws.Range("O2").Value = 0.6283394383394383
$ws.Range("P2").Value = 51.461
$ws.Range("Q2").Value = 0.04234079315451703
$ws.Range("R2").Value = 0.6283394383394383
$ws.Range("U2").Value = 167.1
$ws.Range("V2").Value = 0.1374856014480829
$ws.Range("W2").Value = 0.08330007119091377
$ws.Range("X2").Value = 0.13119530257857
$ws.Range("Y2").Value = -0.04789523138765621
$ws.Range("Z2").Value = 0.04751632577020376
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04240486615532164
$ws.Range("AC2").Value = -0.04240486615532164
$ws.Range("AD2").Value = 5288.9
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 5288.9
$ws.Range("AG2").Value = 5121.799999999999
$ws.Range("AH2").Value = 0.8131390003536123
$ws.Range("AI2").Value = 0.8199717833832034
$ws.Range("AJ2").Value = 0.8082118285678218
$ws.Range("AK2").Value = 0.815183829380869
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

$ws.Range("D3").Value = -0.05980000000000001
$ws.Range("E3").Value = -0.0365
$ws.Range("F3").Value = 0.15
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 51
$ws.Range("L3").Value = 0.3596614950634697
$ws.Range("M3").Value = 51.43
$ws.Range("N3").Value = 0.06182976677085839
$ws.Range("O3").Value = 1.00843137254902
$ws.Range("P3").Value = 51.43
$ws.Range("Q3").Value = 0.06182976677085839
$ws.Range("R3").Value = 1.00843137254902
$ws.Range("U3").Value = 26.1
$ws.Range("V3").Value = 0.03137773503245973
$ws.Range("W3").Value = 0.07731958762886598
$ws.Range("X3").Value = 0.1366022809393718
$ws.Range("Y3").Value = -0.05928269331050581
$ws.Range("Z3").Value = 0.03794690644401628
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.0419249465116964
$ws.Range("AC3").Value = -0.0419249465116964
$ws.Range("AD3").Value = 3748
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 3748
$ws.Range("AG3").Value = 3721.9
$ws.Range("AH3").Value = 0.8183763483121533
$ws.Range("AI3").Value = 0.8317429319604102
$ws.Range("AJ3").Value = 0.8173353536684455
$ws.Range("AK3").Value = 0.8307627061895939
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30.9
$ws.Range("L4").Value = 0.3074626865671642
$ws.Range("M4").Value = 0.031
$ws.Range("N4").Value = 0.0000808133472367049
$ws.Range("O4").Value = 0.001003236245954693
$ws.Range("P4").Value = 0.031
$ws.Range("Q4").Value = 0.0000808133472367049
$ws.Range("R4").Value = 0.001003236245954693
$ws.Range("U4").Value = 141
$ws.Range("V4").Value = 0.367570385818561
$ws.Range("W4").Value = 0.08928055475296157
$ws.Range("X4").Value = 0.1257883242177681
$ws.Range("Y4").Value = -0.03650776946480658
$ws.Range("Z4").Value = 0.07376146788990827
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04288478579894688
$ws.Range("AC4").Value = -0.04288478579894688
$ws.Range("AD4").Value = 1540.9
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1540.9
$ws.Range("AG4").Value = 1399.9
$ws.Range("AH4").Value = 0.800675500129904
$ws.Range("AI4").Value = 0.7926848088893461
$ws.Range("AJ4").Value = 0.7849172974488366
$ws.Range("AK4").Value = 0.7764712407787454
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
